# Daily attendance processing - clear "Recorded By" data in column G
# and shrink the column width now that it no longer needs to fit names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all the recorded-by names in column G (rows 2-177) while keeping
# the cells/styles in place.
$ws.Range("G2:G177").ClearContents()

# Shrink column G from its old width (which accommodated long name lists)
# down to a narrow width. ColumnWidth uses character-width units that get
# internally rounded by the engine, so 12.14 lands on a stored width of 13.
$ws.Columns.Item(7).ColumnWidth = 12.14
